# Generate Report for Handoff
#
# The CI job that produces this localization-status report was rerun.
# That refreshed the "Latest Handoff Date(time)" column for the rows whose
# handoff had not yet settled on a final value, collapsing the previously
# distinct/stale timestamps onto a single new timestamp per sheet:
#
#   Overview!D  (Latest Handoff Date)      -> 2016-03-21 10:26:44
#   zh-cn!E     (Latest Handoff Datetime)  -> 2016-03-21 10:26:40
#   de-de!E     (Latest Handoff Datetime)  -> 2016-03-21 10:26:44
#
# Affected rows (1-based, header is row 1): 7, 10, 11, 12, 13, 14, 15, 16

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-21 10:26:44"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-21 10:26:40"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-21 10:26:44"
}
